$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 values are re-rounded to a coarser ("custom") accuracy of 2 decimals.
$row5 = @{
    "B5"  = 1.62
    "C5"  = 1.27
    "D5"  = 0.01
    "E5"  = 2.69
    "F5"  = 2.45
    "G5"  = 0.79
    "H5"  = 10
    "I5"  = 1.54
    "J5"  = 1.25
    "K5"  = 1.14
    "L5"  = 1.27
    "M5"  = 1.52
    "N5"  = 0.7
    "O5"  = 0.95
    "P5"  = 1.78
    "Q5"  = 0.65
    "R5"  = 0.27
    "S5"  = 0.01
    "T5"  = 9.789999999999999
    "U5"  = 3.31
    "V5"  = 1.24
    "W5"  = 2.56
    "X5"  = 1.02
    "Y5"  = 0.15
    "Z5"  = 4.7
    "AA5" = 0.86
    "AB5" = 0.86
    "AC5" = 1.49
    "AD5" = 1.35
    "AE5" = 0.57
    "AF5" = 9.65
    "AG5" = 0.46
    "AH5" = 1.15
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# The dataset was trimmed (1000 rows across the workbook) - here that drops
# the last data row (row 6), shrinking the used range to A1:AH5.
$ws.Rows("6:6").Delete()
